$d = $word.ActiveDocument

# 1) Add a new row to the end of the first table: date 15/08/2020 with
#    two bullet items ("Mejora ER" and "Mejora prototipos html").
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()

$dateCell = $newRow.Cells.Item(1)
$dateCell.Range.Text = "15/08/2020"

$objCell = $newRow.Cells.Item(2)
$objCell.Range.Text = "Mejora ER" + [char]13 + "Mejora prot"

# Locate the freshly inserted "Mejora prot" text (ranges captured before
# this insertion go stale) and append the remainder of the second bullet.
$tail = $d.Content
$tail.Find.Execute("Mejora prot", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Collapse(0)
$tail.InsertAfter("otipos html")

# 2) Remove the "Mejorar ER" TODO bullet (it now lives in the table above).
$hit = $d.Content
$hit.Find.Execute("Mejorar ER", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $hit.Paragraphs.Item(1)
$pStart = $para.Range.Start
$pEnd = $para.Range.End
$d.Range($pStart, $pEnd + 1).Delete()
